$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 61

# Row 3
$ws.Range("C3").Value = 56

# Row 4
$ws.Range("B4").Value = "<gol>"
$ws.Range("C4").Value = 61

# Row 5
$ws.Range("B5").Value = "<come>"
$ws.Range("C5").Value = 59

# Row 6
$ws.Range("B6").Value = "<loon>"

# Row 7
$ws.Range("B7").Value = "<yes>"

# Row 8
$ws.Range("C8").Value = 58

# Row 10
$ws.Range("B10").Value = "<come>"

# Row 11
$ws.Range("B11").Value = "<way>"

# Row 12
$ws.Range("B12").Value = "<he>"
$ws.Range("C12").Value = 61

# Row 13
$ws.Range("B13").Value = "<the>"
$ws.Range("C13").Value = 63

# Row 14
$ws.Range("B14").Value = "<der>"
